# Laboratorio 9 - Entrega final
# Fills in the student codes and answers to the analysis questions.

$d = $word.ActiveDocument

# --- Student codes: replace the "XXXX" placeholders -------------------
$p1 = $d.Paragraphs(2).Range
$p1.Find.Execute("XXXX", $true, $false, $false, $false, $false, $true, 1, $false, "202013683", 2)

$p2 = $d.Paragraphs(3).Range
$p2.Find.Execute("XXXX", $true, $false, $false, $false, $false, $true, 1, $false, "202013338", 2)

# --- Helper to fill an empty "List Paragraph" answer paragraph --------
function Set-AnswerText($paraIndex, $text) {
    $para = $d.Paragraphs($paraIndex)
    $para.Range.Text = $text
    $para.Range.Font.Name = "Dax-Regular"
    $para.Range.LanguageID = "es-CO"
}

# Question 1 answer (paragraph 6 overall / 1st blank line after Q1)
Set-AnswerText 6 "Si se toma el logaritmo base 2 del numero de elementos del árbol nos dará la altura optima para este árbol. Al hacerlo vemos que el numero obtenido es aproximadamente 10 y se puede ver que esta bastante lejos de 29 (la altura del BST), mientras que la altura del RBT es 13, la cual esta mucho mas cerca de esta altura optima, debido a que el RBT es mas balanceado"

# Question 2 answer (paragraph 9 overall / 1st blank line after Q2)
Set-AnswerText 9 "Seria mayor, ya que así la tabla de hash estuviera ordenada, cada vez que le pida encontrar una fecha, debería pasar por toda la tabla, mientras que el BST, cada vez que busca una fecha, descarta una parte del árbol, lo cual hace que el tiempo sea considerablemente menor."

# Question 3 answer (paragraph 12 overall / 1st blank line after Q3)
Set-AnswerText 12 "Se utiliza la función getCrimesByRange, la cual toma el numero de crímenes del rango y los añade a una lista, luego retorna esa cantidad."

# --- New trailing paragraph with the final remark ----------------------
$endRng = $d.Paragraphs($d.Paragraphs.Count).Range
$endRng.Collapse(0)
$endRng.InsertParagraphAfter()

$newIdx = $d.Paragraphs.Count
$newPara = $d.Paragraphs($newIdx)
$newPara.Style = "Normal"
$newPara.Format.LeftIndent = 18
$newPara.Format.SpaceAfter = 0
$newPara.Alignment = 3

$newPara.Range.Text = "a,b. El RBT tiene menor altura, y esto se debe a que el RBT es capaz de “autobalancearse”, haciendo que los elementos se distribuyan de manera mas equitativa y evitando que se cree mucha mas altura en el arbol de la necesaria."
$newPara.Range.Font.Name = "Dax-Regular"
$newPara.Range.LanguageID = "es-CO"
